# Update crypto price/volume symbol data on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new text value (values are stored as text, not numbers)
$updates = @{
    "D2"  = "247.57"
    "D3"  = "21.80"
    "D5"  = "0.05644"
    "D6"  = "3.389"
    "D7"  = "6.440"
    "D8"  = "0.8013"
    "D9"  = "1.038"
    "D10" = "0.1426"
    "D11" = "0.07226"
    "D12" = "0.03139"
    "D13" = "0.02942"
    "D14" = "0.09284"
    "D15" = "0.001677"
    "D16" = "3.283"
    "D17" = "0.04743"
    "D18" = "0.0005829"
    "D19" = "0.006495"
    "E20" = "19HotbitTokenHTBBestin24h"
    "D21" = "0.001052"
    "D22" = "0.0001503"
    "D27" = "0.1310"
    "D40" = "0.04074"
    "D41" = "0.006915"
    "E41" = "40KickTokenKICK"
    "D42" = "0.1042"
    "D43" = "0.002976"
    "D44" = "0.009124"
    "D47" = "0.7878"
    "D48" = "0.01672"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text interpretation so numeric-looking strings (e.g. "247.57")
    # are not coerced into numbers, then restore the default "Normal"
    # style so no lingering custom number format stays applied to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
